# Auto-generated Excel COM-interop edit script
# Updates the headline/event rows (2-12) on the active worksheet
# to match the new scraped content, expanding the sheet from 9 to 12 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '쓱배송 days'
$ws.Range("B2").Value = '10% 쿠폰 + 카드혜택'
$ws.Range("C2").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003979'
$ws.Range("D2").Value = '1/26~27 쓱배송DAYS'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2023-01-26'
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '2023-01-27'
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = '[''이벤트/쿠폰 > 1/26~27 쓱배송DAYS'', ''스마일클럽'', ''1/26~27 쓱배송DAYS'', ''할인에 할인을 더하다'', ''쓱배송DAYS'', ''최대 10% 장바구니 쿠폰 & 카드할인 바로보기'', ''오늘은 장보는 날 쓱배송 특가상품 바로보기'', ''1. 쓱배송, 새벽배송 모두 쓸 수 있는 10% 장보기 쿠폰'', ''이마트, 트레이더스 쓱배송 및 점포택배/새벽배송'', ''10% 장바구니 쿠폰 - 8만원 이상 결제시 최대 1만 5천원 할인 ID당 1일 1매 (오전 9시 오픈, 선착순 2만명)'', ''10% 장바구니 쿠폰 - 8만원 이상 결제시 최대 2만원 할인 ID당 1일 1매 (오전 9시 오픈, 선착순 5천명, smile club 전용)'', ''스마일클럽 전용 쿠폰 받기'', ''       쿠폰발급 및 사용 기간'', ''       2023년 1월 26일/27일 양일간 (*당일 발급된 쿠폰은 당일 한정 사용가능)'', ''       쿠폰 발급 및 사용일'', ''       23년 1월 26일/27일 양일간 (기간 중 ID당 1일 1회 다운가능, 발급된 쿠폰 당일 한정 사용 가능)'', ''       쿠폰 사용조건'', ''       대상상품 : 이마트몰 쓱배송 / 트레이더스몰 쓱배송 및 점포택배 / 새벽배송 상품'', ''       쿠폰 발급 대상'', ''       스마일클럽 전용 2만원 할인 장바구니 쿠폰 : 스마일클럽 멤버십 가입 회원 전용 (두 쿠폰 모두 이메일주소, SNS 등으로 가입된 간편회원은 제외됩니다)'', ''       쿠폰 적용 제외 상품은 구매금액에 합산되지 않습니다.'', ''+ 결제할 때 한 번 더 할인받기 오늘의 청구할인 카드혜택'', ''[SSGPAY 전용]삼성카드 최대 5만원 5% 청구할인 자세히 보기'', ''3. 오늘은 장보는 날! 쓱배송데이 장보기 특가상품'', ''점포의 신선함을 그대로 간직한 이마트몰 쓱배송 추천상품'', ''가성비 챙기고 대용량 장보기 트레이더스 쓱배송 추천상품'', ''신세계상품권을 SSG MONEY로 전환하고 SSG.COM에서 쇼핑하세요 바로 가기'']'

# Row 3
$ws.Range("A3").Value = '2023 쇼핑 버킷리스트'
$ws.Range("B3").Value = '최대 3만원 할인쿠폰 제공'
$ws.Range("C3").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003961'
$ws.Range("D3").Value = '2023 쇼핑 버킷리스트'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2023-01-25'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '2023-01-31'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = '[''이벤트/쿠폰 > 2023 쇼핑 버킷리스트'', ''스마일클럽'', ''올해 첫 합리적 소비를 위한 장바구니쿠폰 2종'', ''장바구니쿠폰 10% 15만원 이상 최대 3만원'', ''장바구니쿠폰 7% 10만원 이상 최대 1만원'', ''쿠폰은 25일 10시에 오픈됩니다.'', ''발급 후 [My SSG] > [쿠폰] > [보유쿠폰] 에서 확인 가능합니다.'', ''[공식] 윈터 세일 최대 50% & 베스트 상품 제안'', ''★레드하트★한정수량★맨투맨/니트 +)즉시할인쿠폰 바로적용'', ''★SSG x 템퍼★ 덴마크 프리미엄 TEMPUR ~베개/매트리스 ~30% OFF'', ''[자라홈] 이불, 쿠션, 러그 등 : 최대 60%할인'', ''세계 최초 모듈가구 USM, 10%쿠폰 할인'', ''[즉시할인+카드혜택+사은품 4종] 삼성 갤럭시 워치4  모음'', ''SSGX불리 최대 단독 증정 혜택~15%'', ''안티에이징&기미잡티케어 ~50% 7만이상 구매시 신세계상품권 증정! 아이크림/나이트크림/기미잡티앰플/에센스'', ''신년맞이 [SSG X 일동제약 퍼스트랩 ] 10종 스킨&헤어케어 ~80%SALE!'', ''명절 연휴에 추천하는 향수 할인전 / 바이레도, 톰포드, 조말론 외'', ''트와이스 지효PICK  V리프팅 괄사 + 진정크림 UP TO 37% +무료배송'', ''♥시즌오프♥FW 인기 아이템 최대 52% OFF'', ''[에고이스트/랩/플라스틱아일랜드] ~70% OFF! 윈터 클리어런스 외 인기 아이템 쇼핑 위크'', ''[캠브리지] X SSG 아놀 단독 할인 + best 아이템 UP TO 40% OFF'', ''★어그 베스트셀러 ~40% 세일★'', ''[CARLYN]23SS 푸퍼 출시!♥신상&베스트 최대 76% SALE + 추가 쿠폰 혜택까지!'', ''블랙야크/디스커버리 外 인기브랜드 다운/패딩 ~70% OFF'', ''풋웨어 BEST 에어맥스97/데이브레이크 ~44%'', ''[아식스 본사] 최대 82% OFF 다운부터 런닝화까지 베스트템 모음전'', ''투어360 골프화 50%할인'', ''[자라홈] 키즈용품, 장난감, 인형 등 : 최대 60%할인'', ''레고 클리어런스 특가! ★최대 40%할인★ (~2/1일까지)'']'

# Row 4
$ws.Range("A4").Value = '블라썸 위크'
$ws.Range("B4").Value = 'CJ/유한킴벌리 등 장보기 대표브랜드 최대 50%할인, 1+1'
$ws.Range("C4").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004037'
$ws.Range("D4").Value = '(1/26~2/1) 블라썸위크'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2023-01-26'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '2023-02-01'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = '[''이벤트/쿠폰 > (1/26~2/1) 블라썸위크'', ''스마일클럽'', ''CJ/유한킴벌리 등 장보기 대표브랜드 최대 50%할인 및 이벤트'', ''한정수량! 오늘만! 단 하루 타임딜 바로보기'', ''최대 50%할인, 1+1, 사은품증정 블라썸관 행사상품/이벤트 바로보기'', ''한정수량! 오늘만! 단 하루 타임딜'', ''(오늘 단 하루, 1+1) CJ 비비고 특양지곰탕/특설렁탕 700G : 할인전 14,900원 할인후 개당 7,450원~'', ''(오늘 단 하루, 50%) 대상 포도씨유 900ml : 할인전 16,900원 할인후 8,450원'', ''(오늘 단 하루, 1+1) 풀무원 풀무원 모짜렐라 핫도그 5입 : 할인전 8,480원 할인후 개당 4,240원~'', ''(오늘 단 하루, 31%, 1+1) LG생활건강 테크 실내건조 세탁세제3L : 할인전 15,900원 할인후 개당 5,450원'', ''(오늘 단 하루, 한정수량 체럼, 41%) 아모레퍼시픽 (체험팩) 스킨유 이노센트 샤워젤/로션 라보에이치 샴푸 : 할인전 4,900원 할인후 2,900원'', ''(오늘 단 하루, 40%) CJ 헷반 현미/흑미 귀리 곤약밥 150G*12개 : 할인전 35,760원 할인후 21,500원'', ''(오늘 단 하루, 1+1) CJ 백설 파스타 소스 4종 : 할인전 5,500원 할인후 개당 2,790원 ~'', ''(오늘 단 하루, 1+1) 풀무원 들기름을 섞어 바삭바삭 고소하게 구워낸 재래김 80g (5g*16) : 할인전 8,480원 할인후 개당 4,240원~'', ''(오늘 단 하루, 1+1) 풀무원 풀무원 동물복지 수비드 닭가슴살 1150G 3종 : 할인전 3,480원 할인후 개당 1,740원'', ''(오늘 단 하루, 45%) 아모레퍼시픽 LG생활건강 퍼퓸스튜디오 1L 플로라피치 : 할인전 8,900원 할인후 4,900원'', ''(오늘 단 하루, 1+1) 대상 청정원 드레싱 2종 : 할인전 3,600원 할인후 개당 1,800원'', ''(오늘 단 하루, 한정수량, 18%) 매일유업 상하치킨 브리또 6입 : 할인전 20,000원 할인후 16,500원'', ''(오늘 단 하루, 33%) 풀무원 풀무원 행복한 동물복지유정란 15구 780g : 할인전 8,980원 할인후 5,980원'', ''(오늘 단 하루, 30%, 한정수량) LG생활건강 피이조겔 지코의 NO.1보습크립 세트 : 할인전 59,800원 할인후 41,860원'', ''(오늘 단 하루, 60%, 한정수량) 아모레퍼시픽 려 더블이펙터 블랙 샴푸/트리트먼트 543ML : 할인전 39,000원 할인후 15,600원'', ''(오늘 단 하루, 1+1) 우리동네 떡볶이/고소한 유부초밥 330G : 할인전 4,900원 할인후 개당 2,980원~'', ''(오늘 단 하루, 1+1) CJ 비비고 진국육수 닭고기/소고기양지 500g : 할인전 3,080원 할인후 개당 1,540원'', ''(오늘 단 하루, 44%, 한정수량) 매일유업 셀렉스 프로핏 웨ㅐ이프로틴 드링크 초콜릿 2박스 : 할인전 69.900원 할인후 38,900원'', ''(오늘 단 하루, 1+1) LG생활건강 엘라스틴 5X 샴푸/트리트먼트 : 할인전 16,900원 할인후 개당 8,450원'', ''(오늘 단 하루, 1+1) LG생활건강 히말라야 핑크솔트 치약 100g*3_아이스카밍 : 할인전 14,900원 할인후 개당 7,450원'', ''(오늘 단 하루, 1+1) CJ 비비고 사골시래기된장국 460g : 할인전 5,480원 할인후 개당 2,740원'', ''(오늘 단 하루, 1+1) 풀무원 자연은맛있다 건면 4종 : 할인전 5,450원 할인후 개당 2,725원'', ''(오늘 단 하루, 1+1) 풀무원 요리육수 3종 1000ML : 할인전 3,980원 할인후 개당 1,990원'', ''(오늘 단 하루, 2개 구매 시 50%) 유한킴벌리 크리넥스 디럭스미니 카카오230매6입 : 할인전 12,700원 할인후 개당 6,350원'', ''(오늘 단 하루, 62%, 한정수량) LG생활건강 실크테라피 인리처드 액션 60ML X 4개 : 할인전 79,600원 할인후 29,900원'', ''(오늘 단 하루, 1+1) CJ 백설 스파게티소스 640G 2종 : 할인전 4,380원 할인후 개당 2,190원'', ''(오늘 단 하루, 18%, 한정수량) 매일유업 마이카페라떼 마일드로어슈거 220ml 10팩 : 할인전 14,500원 할인후 11,900원'', ''(오늘 단 하루, 25%, 한정수량) 아모레퍼시픽 일리윤 세라마이드 아토로션/울트라리페어로션 + 일리윤 로션 : 할인전 24,500원 할인후 18,743원'', ''(오늘 단 하루, 50%, 한정수량) 유한킴벌리 좋은느낌 유기농순면 내추럴코어 슬날대 14X6 : 할인전 35,900원 할인후 17,950원'', ''(오늘 단 하루, 2개 구매 시, 50%) 유한킴벌리 그린핑거 판테딘 더마 밴드/팬티 기저귀 : 할인전 30,450원 할인후 개당 15,225원~'', ''(오늘 단 하루, 1+1) 풀무원 다논 요거톡125G*4입 : 할인전 3,480원 할인후 개당 1,740원'', ''(오늘 단 하루, 31%, 한정수량) 매일유업 셀렉스 코어프로틴 음료 오리지널 2박스 : 할인전 58,000원 할인후 39,900원'', ''(오늘 단 하루, 2개 구매 시, 50%) 유한킴벌리 크리넥스 안심 프리미엄 3겹 키친타올 140매x4롤 : 할인전 9,200원 할인후 개당 4,600원'', ''(오늘 단 하루, 52%, 한정수량) LG생활건강 세이프 뿌리는 레몬식초 400ml x 3개 : 할인전 22,500원 할인후 10,900원'', ''(오늘 단 하루, 한정수량 체험, 64%) 아모레퍼시픽 롱테이크 샌달우드 인텐시브 3종 KIT : 할인전 10,800원 할인후 3,900원'', ''최대 50%할인 / 1+1 / 다다익선 / 사은품증정'', ''블라썸관 행사상품 및 이벤트'', ''블라썸관 좋아요 누르면 5% 할인쿠폰'', ''최대 50% 할인 행사상품'', ''※ 자세한 내용은 상품 상세페이지 내 이벤트내용을 참고해주세요.'', ''최대 50%할인 행사상품'', ''블라썸관 좋아요 누르면 10% 할인쿠폰'', ''최대 50%할인 행사상품, 기획상품'', ''최대 50%할인, 1+1 행사상품'', ''최대 30%할인, 1+1 행사상품'', ''3만원 이상 구매고객 스타벅스 아메리카노 증정(선착순 1천명)'', ''최대 30%할인, 1+1 행사상품, SSG단독 상품'', ''구매왕 응모이벤트 : 1등 건조기(1명), 2등 스타벅스 아메리카노(100명)'']'

# Row 5
$ws.Range("A5").Value = '더 강력해진 SSG.COM 삼성카드'
$ws.Range("B5").Value = 'SSG MONEY 최대 15% 적립 + 스마일클럽 월이용료 할인'
$ws.Range("C5").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000002385&recruitmentPath=SSG'
$ws.Range("D5").Value = 'SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2022-10-26'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '2025-10-25'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = '[''이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'']'

# Row 6
$ws.Range("A6").Value = '최대 5만원 혜택'
$ws.Range("B6").Value = '+ 쓸 때마다 최대 12% 적립'
$ws.Range("C6").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000000858&siteNo=6005&recruitmentPath=L6007001&eventCode=HPG02'
$ws.Range("D6").Value = 'SSG.COM카드 Edition 2 이벤트 안내 페이지'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2022-07-08'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '2023-06-28'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = '[''이벤트/쿠폰 > SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''SSG.COM카드 Edition2는 SSG.COM에서 최대 5만원 혜택 + 쓸 때마다 최대 12% 적립'', ''이달의 혜택 01. SSG.COM카드 Edition2 첫결제 시 1만원 쿠폰 할인 바로보기'', ''혜택 01. SSG머니 최대 12% 적립 바로보기'', ''혜택 02. 매월 스마일클럽 가입비 3,900원 지원 바로보기'', ''이달의 혜택 01. SSG.COM에서 SSG.COM카드 Edition2 첫 결제 시 1만원 할인'', ''      직전 6개월간 (2022년 7월 1일 부터 2022년 12월 31일) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한함'', ''첫 결제 쿠폰 이용 방법'', ''② SSGPAY에 등록된 SSG.COM카드 Edition2로'', ''③ SSG.COM에서 기간 내 1만 1천원 이상 첫 결제 시 1만원 쿠폰 할인'', ''      본 혜택은 SSG.COM카드 및 SSG.COM카드 Edition2 로 SSG.COM(이마트몰, 신세계몰, 신세계백화점몰 등)에서 직전 6개월간 (2022년 7월 1일 부터 년 12월 31일 까지) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한해 제공됨'', "      본 혜택은 결제 시 [결제방법 > SSGPAY카드] 내 SSG.COM카드 Edition2 선택 시, ''카드할인 최적가'' 추천에 의해 할인 금액이 자동 적용됨. 단, 1만 1천원 이상 첫 결제 시 적용)", ''다운받은 쿠폰은 SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''쿠폰 사용하러 가기'', ''SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''* 쿠폰할인, 적립금 사용 등 선할인을 제외한 카드 최종 결제금액이 2만원 이상인 경우에 한해 혜택 적용'', ''혜택 01. 장 볼 때마다 SSG머니 최대 12% 적립'', ''SSG.COM에서 최대 12%'', ''       스마일클럽으로 5% 적립'', ''       쓱·새벽·트레이더스 구매 시 (구매 전 스마일클럽 적립 아이콘을 꼭 확인해주세요)'', ''       SSG.COM 카드 Edition2로 7% 적립'', ''어디서나 한도없이 0.5%'', ''SSG.COM카드 Edition2로 어디서나 한도없이 0.5% 적립 (SSG.COM 외 모든 가맹점)'', ''SSG Money 최대 12% 적립'', ''스마일클럽 5% 적립 + SSG.COM카드 Edition2 최대 7% 적립'', ''      스마일클럽 5% 적립은 쓱·새벽배송·트레이더스 이용 시에 한함'', ''      SSG.COM카드 Edition2 최대 7% 적립'', ''      SSG.COM에서 결제 시 7%(1만 쓱머니 한도), 그 외 가맹점 0.5% 적립(적립한도 제한 없음)'', ''      무이자 할부 및 현대카드에서 제공하는 다른 할인 서비스 이용 시 적립 제외'', ''스마일클럽 가입비 3,900원 매월 100% 지원'', ''      1. 스마일클럽 자동 가입에'', ''      2. SSGPAY 내 카드 자동 등록'', ''      3. 스마일클럽 정기결제수단 자동 등록 및 월 이용료 3,900원 지원까지! (단, 해당 카드를 월 정기결제 수단에 등록한 경우에 한함)'', ''SSG.COM카드 Edition2를 스마일클럽 월 정기결제 수단에 등록 및 전월 이용금액 30만원 이상 시 혜택 제공'', ''TIP. 스마일클럽 가입 시 SSG.COM 혜택'', ''(SSG 가입 시) 15,000원만 담아도 쓱 · 새벽배송 무료배송'', ''쓱 · 새벽배송 · 트레이더스 SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''스마일클럽 단독 혜택 스타벅스 상품 전용 딜'', ''스마일클럽 가입비 매월 3,900원 지원'', ''월 1회, 매달 스마일클럽 정기결제일에 혜택 제공'', ''SSG.COM카드 Edition2는 최초 발급 시, 스마일클럽 월 정기결제 수단에 자동 등록 됨'', ''전월 이용금액 30만원 미만 시, SSG.COM카드 Edition2로 스마일클럽 정기 결제 금액이 자동 결제됨'', ''스마일클럽 무료 이용 기간이라면 정기결제 금액 지원 대신 SSG머니 3,900원 제공'', ''      스타벅스 자동 충전, 생활요금(통신요금, 아파트관리비 등) 정기결제 신청 및 이체 시 최대 1만원 청구 할인'', ''      2. 스타벅스 자동 충전 또는 생활요금 정기결제 신청(각 항목당 할인한도 5천원, 최대 1만원 할인)'', ''      정기결제 신청 후 카드 결제일에 따라 매출 발생 다음 달 또는 다다음 달 청구 할인 혜택 적용'', ''      단, 청구 할인 제공 일정은 당사 또는 신청인 사정에 의해 상이할 수 있음'', ''      3. 쏘카 1만원 할인쿠폰'', ''      쿠폰은 등록일 포함 30일간 이용 가능'', ''실물 SSG.COM카드 Edition2 수령 전 SSGPAY로 결제 시 건당 100만원 이하 결제 가능 *단, 본인 확인(신분증 확인 및 1원 인증) 완료한 경우에 한하며, 건당 100만원 초과 시 실물카드 수령 후 결제 가능'', ''SSG머니 최대 적립 12%에서 5%는 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능함'', ''카드 이용대금 연체 시 약정금리 + 연체가산금리 3%의 연체이자율이 적용됩니다. (회원별, 이용 상품별 차등적용 / 법정 최고금리 20% 이내) 단, 연체 발생시점에 약정금리가 없는 경우 아래와 같이 적용'', ''일시불 : 거래 발생시점 기준 최소 기간 (2개월)의 유이자 할부 약정금리 + 연체가산금리 3%'', ''무이자할부 : 거래발생시점 기준 동일한 할부 계약 기간의 유이자할부 약정금리 + 연체가산금리 3%'']'

# Row 7
$ws.Range("A7").Value = '쓱배송의 세계로 초대합니다'
$ws.Range("B7").Value = '첫 구매 쿠폰 받고, 친구 초대하면 혜택이 두 배'
$ws.Range("C7").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003900'
$ws.Range("D7").Value = '1월 쓱배송 친구 초대 이벤트'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2023-01-19'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '2023-01-31'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = '[''이벤트/쿠폰 > 1월 쓱배송 친구 초대 이벤트'', ''스마일클럽'', ''1월 쓱배송 친구 초대 이벤트'', ''5천원 할인에 무료배송'', ''친구 초대 이벤트'', ''쓱배송이 처음이라면? 할인에 무료배송'', ''쓱배송 웰컴 쿠폰'', ''5,000원 장바구니 쿠폰 - 첫구매 전용'', ''무료배송 - 첫구매 전용'', ''이마트몰, 쓱배송/새벽배송 20,000원 이상 구매시 사용 가능'', ''첫구매 쿠폰이 모두 소진되었습니다. 더 좋은 혜택으로 찾아뵙겠습니다.'', ''발급 대상 : 2022년 1월 18일 이후 이마트몰, 트레이더스 쓱배송/점포 택배 및 새벽배송 구매 이력이 없는 고객 '', ''쿠폰 사용 조건 : 이마트 쓱배송/새벽배송 상품 2만원 이상 구매 시 사용 가능'', ''이마트몰 쓱배송, 트레이더스 쓱배송, 새벽배송 세가지 배송 모두 한 번도 해본 적 없는 친구에게 나의 초대장 번호를 공유해주세요'', ''응모 이후 초대받은 친구가 이벤트 기간 내 쓱배송 첫 구매를 완료해야만 추천인과 초대받은 사람 모두에게 SSG MONEY 5천원이 적립됩니다. (*첫 구매 완료 기준 - 첫 구매 배송 완료 시점)'', ''SSG MONEY 적립 이후 초대받은 친구가 첫 쓱배송 주문을 전체 반품할 경우, 해당 초대 건에 대하여 SSG MONEY가 회수될 수 있습니다.'', ''본인 스스로 초대한 경우에는 이벤트 참여에서 제외됩니다.'', "친구가 전달 준 초대장 번호를 입력하고 ''첫 쓱배송 응모하기'' 버튼을 눌러주세요 응모 이후, 반드시 이벤트 기간 내 첫 쓱배송 구매를 완료해야만 SSG MONEY가 정상 적립됩니다.", ''첫 쓱배송 응모하기'', ''이벤트 응모여부 확인하기'', ''STEP2 첫 쓱배송 배송 주문하기 : 쓱배송 or 쓱배송TRADERS or 새벽배송(* 주문시, 쓱배송 혹은 새벽배송 스티커를 확인해주세요)'', ''STEP3 SSG MONEY적립 : 이벤트 종료 후 2월 10일(금) 이내 SSG MONEY가 일괄 적립됩니다.'', ''응모 이후 이벤트 기간 내 "쓱배송" 첫 구매를 완료해야만 추천인과 초대받은 사람 모두에게 SSG MONEY 5천원이 적립됩니다. (*첫 구매 완료 기준 - 첫 구매 배송 완료 시점)'', ''쓱배송이 아닌 다른 배송 유형으로 구매했을 시, 참여 인정되지 않습니다. (반드시 쓱배송 혹은 새벽배송 스티커를 확인해주세요)'', ''SSG MONEY 적립 이후 첫 쓱배송 주문을 전체 반품할 경우, 해당 초대 건에 대하여 SSG MONEY가 회수될 수 있습니다.'', ''이벤트 기간 :  1월 19일(목) 00:00 ~ 1월 31일(화) 23:59'', ''최근 1년 간 이마트몰 쓱배송, 트레이더스 쓱배송, 새벽배송 모두 구매 이력이 없는 경우만 쓱배송 첫 구매자로 응모 가능합니다. (2022년 1월 18일 이후 주문 내역이 없는 자)'', ''이벤트 기간 내 등록할 수 있는 초대자는 한 명입니다.'', ''본인 스스로 초대/응모한 경우에는 이벤트 참여에서 제외됩니다.'', ''응모 이후  초대받은 사람이 이벤트 기간 내 쓱배송 첫 구매를 완료해야만 추천인과 초대받은 사람 모두에게 SSG MONEY 5천원이 적립됩니다. (*첫 구매 완료 기준 ? 첫 구매 배송 완료)'', ''SSG MONEY 적립 이후 첫 쓱배송 주문을 전체 반품할 경우, 해당 초대 건에 대하여 추천인과 초대받은 사람 모두 SSG MONEY가 회수될 수 있습니다.'', ''초대 받은 친구가 탈퇴 후 재가입하여 구매를 하더라도 첫 구매로 인정되지 않으며, SSG MONEY가 지급되지 않습니다.'', ''본 이벤트를 통해 지급된 SSG MONEY는 지급일로부터 30일 이후 자동 소멸됩니다.'', ''부정한 방법으로 이벤트에 참여한 것이 발견될 경우, 당첨이 취소될 수 있습니다.'', ''본 이벤트에 대한 세부 사항은 당사 사정에 따라 임의로 변경 혹은 조기 종료될 수 있습니다.'', ''신세계상품권을 SSG MONEY로 전환하고, SSG.COM에서 쇼핑하세요'']'

# Row 8
$ws.Range("A8").Value = '스마일클럽 1월의 가입 혜택'
$ws.Range("B8").Value = 'WELCOME 5천원 할인쿠폰 제공'
$ws.Range("C8").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003449'
$ws.Range("D8").Value = '[스마일클럽] 1월 매일매일 스마일'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2023-01-02'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '2023-01-31'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = '[''이벤트/쿠폰 > [스마일클럽] 1월 매일매일 스마일'', ''스마일클럽'', ''[스마일클럽] 1월 매일매일 스마일'', '' 지금 스마일클럽 가입하고 모든 혜택 받기 '', ''한달 내내 막강한 스마일클럽 혜택'', ''WELCOME 스마일클럽이 처음이라면 첫달 무료 + 5천원 할인 쿠폰 발급 (5,100원 이상 구매 시 사용 가능)'', ''본 쿠폰은 1/5(목)~1/31(화) 기간 내 SSG.COM에서 스마일클럽에 최초 신규 가입한 고객님에 한해 ID당 1회 발급됩니다.'', ''쿠폰은 가입 차주 금요일 이내 자동 지급되며, 앱푸시나 문자메시지 등을 통해 별도 안내 예정입니다. (단, 마케팅 정보 수신을 비동의 하신 경우 안내 발송이 제한됩니다.)'', ''본 쿠폰은 5,100원 이상 구매 시 5,000원 할인됩니다. (할인액 및 배송비 제외한 구매 금액 기준으로 쿠폰 적용됨)'', ''일부 상품 및 브랜드는 쿠폰 적용 제외될 수 있습니다.'', ''SSG.COM의 쿠폰은 결제 시 최적 할인에 의해 자동 설정되며, 직접 변경이 가능합니다.'', ''본 이벤트는 당사 사정으로 내용이 변경되거나 종료될 수 있습니다.'', ''SSG.COM 카드 Edition2 매월 스마일클럽 가입비 100% 지원'', ''장 보는 날에도 멤버십은 더 큰 혜택 쓱배송데이 최대 2만원 10% 할인 쿠폰'', ''2023년 검은 토끼의 해 맞이 건강 프로젝트 스마일클럽 전용 최대 15% 쿠폰으로 건강을 선물해요'', ''계속해서 이어질 스마일클럽 전용 혜택을 기대해주세요'']'

# Row 9
$ws.Range("A9").Value = '쓱스럽게 안녕'
$ws.Range("B9").Value = 'See you again, Say Goodbye'
$ws.Range("C9").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003967'
$ws.Range("D9").Value = '쓱스럽게 안녕'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2023-01-17'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = '2023-01-31'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = '[''이벤트/쿠폰 > 쓱스럽게 안녕'', ''스마일클럽'']'

# Row 10
$ws.Range("A10").Value = '민지의 장바구니-겨울 먹거리편'
$ws.Range("B10").Value = '디저트 레시피 영상보고 댓글 달면 SSG머니 추첨 증정'
$ws.Range("C10").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003980'
$ws.Range("D10").Value = '민지의 장바구니-겨울 먹거리편 (1/26~2/1)'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2023-01-26'
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = '2023-02-01'
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = '[''이벤트/쿠폰 > 민지의 장바구니-겨울 먹거리편 (1/26~2/1)'', ''스마일클럽'', ''SSG MONEY 1천원 (500명)'', ''       My SSG > 이벤트 현황 > 이벤트 참여내역에서 확인해주세요'', ''       - ID당 이벤트 기간 내 1회 응모 가능하며, 당첨은 기간 내 ID 당 1회만 가능합니다.'', ''       - 이벤트 성격과 다른 내용이거나 불쾌감을 줄 수 있는 댓글은 사전 고지 없이 삭제될 수 있습니다.'', ''최대 2만원 할인 10% 장보기 쿠폰 받고 오세요!'', ''쓱배송데이 보러가기'', ''신세계몰 상품권을 SSG MONEY로 전환하고 SSG.COM에서 쇼핑하세요'']'

# Row 11
$ws.Range("A11").Value = '신일전자 1/26(목) 8PM'
$ws.Range("B11").Value = '신일 팬히터 1200 아이보리 선착순 핫딜 39만원대 + 전고객 캠핑용 전기그릴 증정'
$ws.Range("C11").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003983'
$ws.Range("D11").Value = '신일 팬히터 선착순 핫딜 39만원대 @SSG.LIVE 1/26(목) 8PM'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2023-01-19'
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = '2023-01-26'
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = '[''이벤트/쿠폰 > 신일 팬히터 선착순 핫딜 39만원대 @SSG.LIVE 1/26(목) 8PM'', ''스마일클럽'', ''신일 팬히터 선착순 핫딜 39만원대 @SSG.LIVE 1/26(목) 8PM'']'

# Row 12
$ws.Range("A12").Value = '피코크 1/27(금) 11AM'
$ws.Range("B12").Value = '요리가 제일 귀찮을 명절 직후, 피코크 간식/간편식 최대 52%'
$ws.Range("C12").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000003986'
$ws.Range("D12").Value = '피코크 @SSG.LIVE 1/27(금) 11:00'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2023-01-17'
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = '2023-01-27'
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = '[''이벤트/쿠폰 > 피코크 @SSG.LIVE 1/27(금) 11:00'', ''스마일클럽'']'

